$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.668.75'
$ws.Range("E2").Value = '  +2.90%  '
$ws.Range("D3").Value = '3.200.33'
$ws.Range("E3").Value = '  +1.72%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'600.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.88%  '
$ws.Range("D6").Value = "'157.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.92%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '3.203.42'
$ws.Range("E8").Value = '  +1.83%  '
$ws.Range("D9").Value = "'0.552"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.57%  '
$ws.Range("E10").Value = '  +1.42%  '
$ws.Range("D11").Value = "'5.99"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.84%  '
$ws.Range("D12").Value = "'0.514"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.36%  '
$ws.Range("D13").Value = "'0.0000267"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("D14").Value = "'39.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.36%  '
$ws.Range("D15").Value = '3.732.91'
$ws.Range("E15").Value = '  +1.99%  '
$ws.Range("D16").Value = '66.750.63'
$ws.Range("E16").Value = '  +2.94%  '
$ws.Range("D17").Value = "'7.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.31%  '
$ws.Range("D18").Value = '3.208.24'
$ws.Range("E18").Value = '  +2.28%  '
$ws.Range("E19").Value = '  +0.80%  '
$ws.Range("D20").Value = "'517.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.55%  '
$ws.Range("D21").Value = "'15.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("D22").Value = "'0.740"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.68%  '
$ws.Range("D23").Value = "'8.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.16%  '
$ws.Range("D24").Value = "'14.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = "'85.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.98%  '
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("D27").Value = "'9.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.90%  '
$ws.Range("E28").Value = '  +3.39%  '
$ws.Range("D29").Value = "'2.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +10.73%  '
$ws.Range("D30").Value = "'3.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.14%  '
$ws.Range("D31").Value = "'7.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.27%  '
$ws.Range("D32").Value = "'28.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.16%  '
$ws.Range("D33").Value = "'1.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("D35").Value = "'6.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.11%  '
$ws.Range("D36").Value = "'525.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +12.14%  '
$ws.Range("D37").Value = "'54.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("D38").Value = "'0.0901"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.87%  '
$ws.Range("D39").Value = "'0.0424"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.62%  '
$ws.Range("E40").Value = '  +9.32%  '
$ws.Range("D41").Value = "'2.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.33%  '
$ws.Range("D42").Value = "'8.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.06%  '
$ws.Range("D43").Value = "'0.307"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.83%  '
$ws.Range("D44").Value = '0.0₃0689'
$ws.Range("E44").Value = '  +15.25%  '
$ws.Range("E45").Value = '  +3.24%  '
$ws.Range("D46").Value = '2.883.08'
$ws.Range("E46").Value = '  -3.37%  '
$ws.Range("D47").Value = "'28.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.90%  '
$ws.Range("D48").Value = "'2.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.20%  '
$ws.Range("D49").Value = "'0.117"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.78%  '
$ws.Range("E51").Value = '  +9.64%  '
